$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (Changed) column C for all existing data rows (2-504)
# from 45188 to 45189.
$ws.Range("C2:C504").Value = 45189

# Row 504 is no longer the last row once new rows are appended below it, so it
# now also gets the explicit row height that every non-final row carries.
$ws.Rows.Item(504).RowHeight = 15

# 2. Append new data rows 505-511.
$newRows = @(
    @{ Row = 505; A = "A 43883-2023"; B = 45187; C = 45189; G = 3.7 },
    @{ Row = 506; A = "A 43897-2023"; B = 45187; C = 45189; G = 1.8 },
    @{ Row = 507; A = "A 44040-2023"; B = 45188; C = 45189; G = 1   },
    @{ Row = 508; A = "A 44086-2023"; B = 45188; C = 45189; G = 2.1 },
    @{ Row = 509; A = "A 44133-2023"; B = 45188; C = 45189; G = 5.1 },
    @{ Row = 510; A = "A 44090-2023"; B = 45188; C = 45189; G = 0.5 },
    @{ Row = 511; A = "A 44087-2023"; B = 45188; C = 45189; G = 1.1 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = "HALLANDS LÄN"
    $ws.Cells.Item($row, 5).Value = "FALKENBERG"

    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0
    $ws.Cells.Item($row, 15).Value = 0
    $ws.Cells.Item($row, 16).Value = 0
    $ws.Cells.Item($row, 17).Value = 0

    $ws.Cells.Item($row, 18).WrapText = $true

    # Rows 505-510 get an explicit row height (matches the rest of the sheet);
    # the final row (511) is left without it, matching the source data.
    if ($row -lt 511) {
        $ws.Rows.Item($row).RowHeight = 15
    }
}
